# Update countries & provincias Spain
# Applies the 27-May-2020 18:35 -> 19:05 data refresh to the "Pais" sheet:
#   - bumps the "Datos actualizados..." timestamp cell
#   - updates the numeric counters for several countries whose figures changed
#   - re-sorts a handful of tied/near-tied low-count countries, which moves
#     their names (and associated row data) to different rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 19:05"

# --- Updated country statistics (Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
Set-Row 4 @("Estados Unidos", 1733961, 8686, 481988, 1150832, 0, 569, 101141)

# India
Set-Row 13 @("India", 157777, 6984, 67452, 85797, 0, 184, 4528)

# Canada
Set-Row 16 @("Canada", 87481, 834, 46072, 34649, 0, 121, 6760)

# Irak
Set-Row 69 @("Irak", 5135, 287, 2904, 2056, 0, 6, 175)

# Mali
Set-Row 108 @("Mali", 1116, 39, 632, 414, 0, 0, 70)

# Jordania
Set-Row 128 @("Jordania", 720, 2, 586, 125, 0, 0, 9)

# Aruba
Set-Row 171 @("Aruba", 101, 0, 97, 1, 0, 0, 3)

# --- Re-sorted low-count tail: rows 199-212 ---
Set-Row 199 @("Belice", 18, 0, 16, 0, 0, 0, 2)
Set-Row 200 @("Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)
Set-Row 201 @("Santa Lucia", 18, 0, 18, 0, 0, 0, 0)
Set-Row 207 @("Islas Turcas y Caicos", 12, 0, 10, 1, 0, 0, 1)
Set-Row 208 @("Groenlandia", 12, 0, 11, 1, 0, 0, 0)
Set-Row 209 @("Surinam", 11, 0, 9, 1, 0, 0, 1)
Set-Row 210 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)
Set-Row 211 @("Montserrat", 11, 0, 10, 0, 0, 0, 0)
Set-Row 212 @("Sahara Occidental", 9, 0, 6, 2, 0, 0, 0)
